# 1 new TC committed.
# Adds a new "noplates" test-case row (row 3) to both the QA and FT sheets,
# each with an email hyperlink in column A and a "Tester@123" hyperlink in
# column B (mirroring the existing row 2 pattern), widens column A on the
# FT sheet, and leaves the QA sheet as the active tab/selection.

$wb = $excel.ActiveWorkbook

$newEmail = "ayush.kumar+noplates@revivertest.com"
$password = "Tester@123"

$qa = $wb.Worksheets.Item(1)
$ft = $wb.Worksheets.Item(2)

# --- QA sheet (sheet1): add row 3 ---
$qa.Range("A3").Value = $newEmail
$qa.Range("B3").Value = $password
$qa.Hyperlinks.Add($qa.Range("A3"), "mailto:" + $newEmail)
$qa.Hyperlinks.Add($qa.Range("B3"), "mailto:" + $password)
$qa.Range("A3").Style = "Hyperlink"
$qa.Range("B3").Style = "Hyperlink"

# --- FT sheet (sheet2): add row 3 ---
$ft.Range("A3").Value = $newEmail
$ft.Range("B3").Value = $password
$ft.Hyperlinks.Add($ft.Range("A3"), "mailto:" + $newEmail)
$ft.Hyperlinks.Add($ft.Range("B3"), "mailto:" + $password)
$ft.Range("A3").Style = "Hyperlink"
$ft.Range("B3").Style = "Hyperlink"

# Widen column A on the FT sheet (matches the author's resize after adding
# the longer "noplates" email address).
$ft.Columns.Item(1).ColumnWidth = 38.330729166666664

# --- Selection / active-tab bookkeeping ---
# Select the new row's range on FT first (mirrors the author highlighting
# the pasted range there) ...
[void]$ft.Range("A3:B3").Select()
# ... then switch to the QA sheet and land on F12, which becomes the
# workbook's active tab/selection at save time.
$qa.Activate()
[void]$qa.Range("F12").Select()
